$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the period date range from 10/06/2025 to 11/06/2025 for all data rows (rows 2-19)
$ws.Range("G2:G19").Value = "11/06/2025 00:00"
$ws.Range("H2:H19").Value = "11/06/2025 23:59"

# Update recalculated metric values for the affected equipment rows
$ws.Range("I2").Value = 19.400000000000546
$ws.Range("J2").Value = 8088.7
$ws.Range("K2").Value = 63.75316611111114
$ws.Range("L2").Value = 225.91780083333344
$ws.Range("M2").Value = 3.4775955555555553
$ws.Range("N2").Value = 6.652043398232149
$ws.Range("O2").Value = 20.35943581424459
$ws.Range("P2").Value = 35.10527265077139
$ws.Range("Q2").Value = 293.1485652777778
$ws.Range("R2").Value = 967.2467659565131
$ws.Range("S2").Value = 0.68
$ws.Range("T2").Value = 14.123999431917145
$ws.Range("U2").Value = 46.60228427162318
$ws.Range("V2").Value = 0.27723785433317394
$ws.Range("I4").Value = 19.350000000000364
$ws.Range("J4").Value = 8121.9
$ws.Range("K4").Value = 42.92891833333333
$ws.Range("L4").Value = 231.72235194444448
$ws.Range("N4").Value = 4.853603143051577
$ws.Range("O4").Value = 19.497712057964243
$ws.Range("Q4").Value = 274.6512716666666
$ws.Range("R4").Value = 906.2147516693477
$ws.Range("S4").Value = 0.3
$ws.Range("T4").Value = 13.277337253372947
$ws.Range("U4").Value = 43.80871353291403
$ws.Range("V4").Value = 0.2752595992742738
$ws.Range("I7").Value = 20.400000000001455
$ws.Range("J7").Value = 8338.45
$ws.Range("K7").Value = 67.70883583333332
$ws.Range("L7").Value = 181.30918999999994
$ws.Range("M7").Value = 0.02260527777777778
$ws.Range("N7").Value = 5.975965089363547
$ws.Range("O7").Value = 17.840973894813946
$ws.Range("P7").Value = 38.269
$ws.Range("Q7").Value = 249.04063527777774
$ws.Range("R7").Value = 821.7121882753635
$ws.Range("S7").Value = 0.608
$ws.Range("T7").Value = 11.600214124440885
$ws.Range("U7").Value = 38.275028177731514
$ws.Range("V7").Value = 0.2839910145897543
$ws.Range("I8").Value = 21.25
$ws.Range("J8").Value = 6348.7
$ws.Range("K8").Value = 55.304571666666654
$ws.Range("L8").Value = 193.68677277777772
$ws.Range("M8").Value = 0.5724213888888889
$ws.Range("N8").Value = 4.950021491339216
$ws.Range("O8").Value = 18.10333113347889
$ws.Range("P8").Value = 27.117877146631436
$ws.Range("Q8").Value = 249.56376138888902
$ws.Range("R8").Value = 823.4382483660221
$ws.Range("S8").Value = 0.52
$ws.Range("T8").Value = 11.495198867374963
$ws.Range("U8").Value = 37.928529235542065
$ws.Range("V8").Value = 0.27292770816897716
$ws.Range("I10").Value = 17.850000000000364
$ws.Range("J10").Value = 8828.95
$ws.Range("K10").Value = 45.80332444444445
$ws.Range("L10").Value = 149.89100833333333
$ws.Range("M10").Value = 0.08514583333333332
$ws.Range("N10").Value = 4.738515685037058
$ws.Range("O10").Value = 17.11427156123925
$ws.Range("P10").Value = 20.291
$ws.Range("Q10").Value = 195.77948083333334
$ws.Range("R10").Value = 645.9764505319987
$ws.Range("S10").Value = 0.7
$ws.Range("T10").Value = 10.688962707586514
$ws.Range("U10").Value = 35.268344569744244
$ws.Range("V10").Value = 0.2610135727322735
$ws.Range("I12").Value = 18.25
$ws.Range("J12").Value = 920.75
$ws.Range("K12").Value = 95.89975750000005
$ws.Range("L12").Value = 177.95285222222222
$ws.Range("M12").Value = 0.3321738888888889
$ws.Range("N12").Value = 10.240484505078964
$ws.Range("O12").Value = 18.650448700669237
$ws.Range("P12").Value = 25.156421276595744
$ws.Range("Q12").Value = 274.1847836111112
$ws.Range("R12").Value = 904.6755694370698
$ws.Range("S12").Value = 0.496
$ws.Range("T12").Value = 14.492057236981184
$ws.Range("U12").Value = 47.81669485997421
$ws.Range("V12").Value = 0.3580057076203627
$ws.Range("I13").Value = 20.299999999999955
$ws.Range("J13").Value = 919.35
$ws.Range("K13").Value = 65.34161027777778
$ws.Range("L13").Value = 222.7585213888889
$ws.Range("M13").Value = 0.10196666666666666
$ws.Range("N13").Value = 7.432575328876699
$ws.Range("O13").Value = 19.3091795331099
$ws.Range("P13").Value = 22.215
$ws.Range("Q13").Value = 288.20209694444446
$ws.Range("R13").Value = 950.9258418073897
$ws.Range("S13").Value = 0.428
$ws.Range("T13").Value = 14.174492627431912
$ws.Range("U13").Value = 46.76888709984508
$ws.Range("V13").Value = 0.32642344539726603
$ws.Range("I19").Value = 16.5
$ws.Range("J19").Value = 928.65
$ws.Range("K19").Value = 93.16018305555555
$ws.Range("L19").Value = 190.67577916666676
$ws.Range("M19").Value = 4.571635
$ws.Range("N19").Value = 10.454452363439158
$ws.Range("O19").Value = 23.64379718820174
$ws.Range("P19").Value = 31.25084423735261
$ws.Range("Q19").Value = 288.4075988888889
$ws.Range("R19").Value = 951.6038976285851
$ws.Range("S19").Value = 0.49200000000000005
$ws.Range("T19").Value = 16.84425112417512
$ws.Range("U19").Value = 55.57778326282945
$ws.Range("V19").Value = 0.38509395188315826
